$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K (before the existing "DC Timestamp" column)
$ws.Columns.Item(11).Insert()

# Set header for the new column
$ws.Cells.Item(1, 11).Value = "Imaging Start"

# Give the new column a "best fit" style width, matching how Excel auto-sizes
# a column to fit its header text (closest representable width in this runtime).
$ws.Columns.Item(11).ColumnWidth = 11.7

# Match new active cell selection as seen in the edited workbook
$ws.Range("K2").Select()
